$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stored "GameExcelPath" value (B2) to the new local file path.
$ws.Range("B2").Value = "C:\Users\david\Desktop\RPA\proiect\Steam-game-finder\SteamGameFinder\games_to_find1.xlsx"

# Leave the active selection on B8, matching the saved view state.
$ws.Range("B8").Select() | Out-Null
